# Update cryptos list with refreshed market data (GitHub Actions bot)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.585.10'
$ws.Range("E2").Value = '  +0.11%  '
$ws.Range("D3").Value = '1.809.20'
$ws.Range("E3").Value = '  -0.10%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("E5").Value = '  -1.35%  '
$ws.Range("E6").Value = '  +2.78%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '37.42'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.94%  '
$ws.Range("E9").Value = '  -3.37%  '
$ws.Range("E10").Value = '  -2.19%  '
$ws.Range("E11").Value = '  +1.48%  '
$ws.Range("D12").Value = '2.070.57'
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.31'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.86%  '
$ws.Range("D14").Value = '1.838.38'
$ws.Range("E14").Value = '  +1.69%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.633'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.05%  '
$ws.Range("D16").Value = '34.544.23'
$ws.Range("E16").Value = '  +0.08%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.43'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.65'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.79%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.42'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.10%  '
$ws.Range("D20").Value = '0.0₃0776'
$ws.Range("E20").Value = '  -3.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.21'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.41%  '
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("E23").Value = '  -0.95%  '
$ws.Range("E24").Value = '  +4.61%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '172.32'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.69%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.86'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.37%  '
$ws.Range("E27").Value = '  +2.65%  '
$ws.Range("E28").Value = '  +1.33%  '
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.82'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.94%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.93'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.92%  '
$ws.Range("E32").Value = '  -1.54%  '
$ws.Range("E33").Value = '  -2.82%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.81'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.62%  '
$ws.Range("D35").Value = '1.365.70'
$ws.Range("E35").Value = '  -2.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.653'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.58%  '
$ws.Range("E37").Value = '  +0.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.37'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.38%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0187'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.84%  '
$ws.Range("E40").Value = '  +1.72%  '
$ws.Range("E41").Value = '  -1.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '80.85'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.939'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.33%  '
$ws.Range("E44").Value = '  +5.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.73'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.30%  '
$ws.Range("E46").Value = '  -1.99%  '
$ws.Range("D47").Value = '1.970.84'
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.81'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.07%  '
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '102.82'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.10%  '
$ws.Range("D51").Value = '0.0₆0122'
$ws.Range("E51").Value = '  -6.29%  '
